{"js": "// Replace the division-problem text in the worksheet table.\n// Each mapping is an exact, unique, whole-cell string, so a simple\n// search-and-replace (old text -> new text) reproduces the diff.\nconst replacements = [\n  [\"43\u00f77=\", \"35\u00f76=\"],\n  [\"63\u00f76=\", \"23\u00f74=\"],\n  [\"60\u00f73=\", \"67\u00f73=\"],\n  [\"35\u00f74=\", \"16\u00f72=\"],\n  [\"37\u00f79=\", \"50\u00f73=\"],\n  [\"33\u00f74=\", \"32\u00f73=\"],\n  [\"86\u00f76=\", \"68\u00f77=\"],\n  [\"92\u00f75=\", \"75\u00f75=\"],\n  [\"57\u00f72=\", \"10\u00f75=\"],\n  [\"78\u00f75=\", \"44\u00f79=\"],\n  [\"22\u00f74=\", \"75\u00f77=\"],\n  [\"63\u00f75=\", \"53\u00f73=\"],\n  [\"64\u00f77=\", \"92\u00f75=\"],\n  [\"46\u00f74=\", \"91\u00f75=\"],\n  [\"99\u00f73=\", \"71\u00f79=\"],\n  [\"10\u00f72=\", \"56\u00f73=\"],\n  [\"63\u00f78=\", \"72\u00f73=\"],\n  [\"15\u00f75=\", \"81\u00f73=\"],\n  [\"80\u00f74=\", \"15\u00f74=\"],\n  [\"16\u00f78=\", \"33\u00f78=\"],\n  [\"97\u00f75=\", \"44\u00f74=\"],\n  [\"36\u00f77=\", \"23\u00f77=\"],\n  [\"48\u00f77=\", \"32\u00f77=\"],\n  [\"84\u00f76=\", \"55\u00f77=\"],\n  [\"50\u00f75=\", \"90\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in the worksheet table.\n# Each mapping is an exact, unique, whole-cell string, so Find/Replace\n# (old text -> new text) across the whole document reproduces the diff.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"43\u00f77=\", \"35\u00f76=\"),\n    @(\"63\u00f76=\", \"23\u00f74=\"),\n    @(\"60\u00f73=\", \"67\u00f73=\"),\n    @(\"35\u00f74=\", \"16\u00f72=\"),\n    @(\"37\u00f79=\", \"50\u00f73=\"),\n    @(\"33\u00f74=\", \"32\u00f73=\"),\n    @(\"86\u00f76=\", \"68\u00f77=\"),\n    @(\"92\u00f75=\", \"75\u00f75=\"),\n    @(\"57\u00f72=\", \"10\u00f75=\"),\n    @(\"78\u00f75=\", \"44\u00f79=\"),\n    @(\"22\u00f74=\", \"75\u00f77=\"),\n    @(\"63\u00f75=\", \"53\u00f73=\"),\n    @(\"64\u00f77=\", \"92\u00f75=\"),\n    @(\"46\u00f74=\", \"91\u00f75=\"),\n    @(\"99\u00f73=\", \"71\u00f79=\"),\n    @(\"10\u00f72=\", \"56\u00f73=\"),\n    @(\"63\u00f78=\", \"72\u00f73=\"),\n    @(\"15\u00f75=\", \"81\u00f73=\"),\n    @(\"80\u00f74=\", \"15\u00f74=\"),\n    @(\"16\u00f78=\", \"33\u00f78=\"),\n    @(\"97\u00f75=\", \"44\u00f74=\"),\n    @(\"36\u00f77=\", \"23\u00f77=\"),\n    @(\"48\u00f77=\", \"32\u00f77=\"),\n    @(\"84\u00f76=\", \"55\u00f77=\"),\n    @(\"50\u00f75=\", \"90\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
